$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Extend header row H1:M1 with the same style as the existing header cells (copy format from A1)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("H1:M1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 2) Cells that hold numeric-looking text must stay text (not auto-converted to numbers),
#    so pre-format them as Text before assigning the value (matches source data being all strings).
$textCells = @("C2","D2","E2","G2","I2","J2","L2","M2","C3","D3","E3","G3","J3","L3","M3","C4","D4","E4","G4","J4","L4","M4","J5","L5","M5","J6","L6","M6","J7","L7","M7")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# 3) Header row updates (B1:M1)
$ws.Range("B1").Value = "test1.csv-COL2"
$ws.Range("C1").Value = "test1.csv-COL3"
$ws.Range("D1").Value = "test1.csv-COL4"
$ws.Range("E1").Value = "http://dbpedia.org/ontology/percentage"
$ws.Range("F1").Value = "test1.csv-COL6"
$ws.Range("G1").Value = "http://dbpedia.org/ontology/populationTotal"
$ws.Range("H1").Value = "test2.csv-COL2"
$ws.Range("I1").Value = "test2.csv-COL3"
$ws.Range("J1").Value = "http://dbpedia.org/ontology/weight"
$ws.Range("K1").Value = "test3.csv-COL3"
$ws.Range("L1").Value = "test3.csv-COL4"
$ws.Range("M1").Value = "test3.csv-COL5"

# 4) Data rows 2-7 updates
$ws.Range("A2").Value = "http://dbpedia.org/resource/Terp"
$ws.Range("B2").Value = "http://dbpedia.org/resource/NAMUR"
$ws.Range("C2").Value = "92136"
$ws.Range("D2").Value = "5100"
$ws.Range("E2").Value = "54"
$ws.Range("F2").Value = "http://dbpedia.org/resource/WD"
$ws.Range("G2").Value = "1845"
$ws.Range("H2").Value = "http://dbpedia.org/resource/A1"
$ws.Range("I2").Value = "92094046"
$ws.Range("J2").Value = "92044"
$ws.Range("K2").Value = "http://dbpedia.org/resource/Terp"
$ws.Range("L2").Value = "5023"
$ws.Range("M2").Value = "4.605"
$ws.Range("A3").Value = "http://dbpedia.org/resource/Terp"
$ws.Range("B3").Value = "http://dbpedia.org/resource/NAMUR"
$ws.Range("C3").Value = "92136"
$ws.Range("D3").Value = "5100"
$ws.Range("E3").Value = "54"
$ws.Range("F3").Value = "http://dbpedia.org/resource/WD"
$ws.Range("G3").Value = "1845"
$ws.Range("J3").Value = "92044"
$ws.Range("K3").Value = "http://dbpedia.org/resource/Terp"
$ws.Range("L3").Value = "5023"
$ws.Range("M3").Value = "4.605"
$ws.Range("A4").Value = "http://dbpedia.org/resource/Flawinne"
$ws.Range("B4").Value = "http://dbpedia.org/resource/NAMUR"
$ws.Range("C4").Value = "92043"
$ws.Range("D4").Value = "5020"
$ws.Range("E4").Value = "71"
$ws.Range("F4").Value = "http://dbpedia.org/resource/FW"
$ws.Range("G4").Value = "4491"
$ws.Range("J4").Value = "92043"
$ws.Range("K4").Value = "http://dbpedia.org/resource/Flawinne"
$ws.Range("L4").Value = "5020"
$ws.Range("M4").Value = "6.742"
$ws.Range("A5").Value = "http://dbpedia.org/resource/Citadelle"
$ws.Range("B5").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("J5").Value = "92075"
$ws.Range("K5").Value = "http://dbpedia.org/resource/Citadelle"
$ws.Range("L5").Value = "5101"
$ws.Range("M5").Value = "3.315"
$ws.Range("A6").Value = "http://dbpedia.org/resource/Flawinne"
$ws.Range("J6").Value = "92043"
$ws.Range("K6").Value = "http://dbpedia.org/resource/Flawinne"
$ws.Range("L6").Value = "5020"
$ws.Range("M6").Value = "6.742"
$ws.Range("A7").Value = "http://dbpedia.org/resource/Terp"
$ws.Range("J7").Value = "92044"
$ws.Range("K7").Value = "http://dbpedia.org/resource/Terp"
$ws.Range("L7").Value = "5023"
$ws.Range("M7").Value = "4.605"
